$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit "last update for day": update the first cell's text from
# "This is" to "weiner". All other cell values (A2:A7) are unchanged.
$ws.Range("A1").Value = "weiner"
